$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $escaped = $val.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue $ws 'D2' '51.479.58'
Set-TextValue $ws 'E2' '  +0.42%  '
Set-TextValue $ws 'D3' '2.980.12'
Set-TextValue $ws 'E3' '  +2.09%  '
Set-TextValue $ws 'E4' '  -0.01%  '
Set-TextValue $ws 'D5' '381.13'
Set-TextValue $ws 'E5' '  +2.44%  '
Set-TextValue $ws 'D6' '105.43'
Set-TextValue $ws 'E6' '  +1.77%  '
Set-TextValue $ws 'E7' '  +0.56%  '
Set-TextValue $ws 'D8' '1.00'
Set-TextValue $ws 'E8' '  +0.00%  '
Set-TextValue $ws 'E9' '  +0.91%  '
Set-TextValue $ws 'D10' '37.25'
Set-TextValue $ws 'E10' '  +0.32%  '
Set-TextValue $ws 'E11' '  +0.47%  '
Set-TextValue $ws 'E12' '  +0.73%  '
Set-TextValue $ws 'D13' '3.448.47'
Set-TextValue $ws 'E13' '  +2.07%  '
Set-TextValue $ws 'E14' '  +0.52%  '
Set-TextValue $ws 'D15' '7.50'
Set-TextValue $ws 'E15' '  +1.70%  '
Set-TextValue $ws 'D16' '2.968.82'
Set-TextValue $ws 'E16' '  +1.92%  '
Set-TextValue $ws 'D17' '0.972'
Set-TextValue $ws 'E17' '  +1.89%  '
Set-TextValue $ws 'D18' '51.478.29'
Set-TextValue $ws 'E18' '  +0.46%  '
Set-TextValue $ws 'E19' '  +0.54%  '
Set-TextValue $ws 'D20' '7.40'
Set-TextValue $ws 'E20' '  +1.38%  '
Set-TextValue $ws 'E21' '  -0.72%  '
Set-TextValue $ws 'E22' '  +1.65%  '
Set-TextValue $ws 'D23' '69.14'
Set-TextValue $ws 'E23' '  +1.03%  '
Set-TextValue $ws 'D24' '262.76'
Set-TextValue $ws 'E24' '  +0.67%  '
Set-TextValue $ws 'E25' '  +3.86%  '
Set-TextValue $ws 'E26' '  -2.40%  '
Set-TextValue $ws 'D27' '7.22'
Set-TextValue $ws 'E27' '  +16.76%  '
Set-TextValue $ws 'D28' '7.51'
Set-TextValue $ws 'E28' '  +1.14%  '
Set-TextValue $ws 'B29' 'Hedera'
Set-TextValue $ws 'C29' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D29' '0.114'
Set-TextValue $ws 'E29' '  +9.79%  '
Set-TextValue $ws 'B30' 'Dai'
Set-TextValue $ws 'C30' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D30' '1.00'
Set-TextValue $ws 'E30' '  +0.03%  '
Set-TextValue $ws 'D31' '26.01'
Set-TextValue $ws 'E31' '  +0.19%  '
Set-TextValue $ws 'D32' '9.86'
Set-TextValue $ws 'E32' '  -1.00%  '
Set-TextValue $ws 'D33' '34.93'
Set-TextValue $ws 'E33' '  -1.14%  '
Set-TextValue $ws 'B34' 'OKB'
Set-TextValue $ws 'C34' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D34' '51.30'
Set-TextValue $ws 'E34' '  +1.03%  '
Set-TextValue $ws 'B35' 'Toncoin'
Set-TextValue $ws 'C35' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws 'D35' '2.09'
Set-TextValue $ws 'E35' '  -2.42%  '
Set-TextValue $ws 'D36' '0.0451'
Set-TextValue $ws 'E36' '  +6.31%  '
Set-TextValue $ws 'E37' '  +0.07%  '
Set-TextValue $ws 'E38' '  -1.38%  '
Set-TextValue $ws 'E39' '  +2.52%  '
Set-TextValue $ws 'E40' '  -5.82%  '
Set-TextValue $ws 'E41' '  -0.84%  '
Set-TextValue $ws 'E42' '  +2.28%  '
Set-TextValue $ws 'D43' '124.70'
Set-TextValue $ws 'E43' '  +6.34%  '
Set-TextValue $ws 'D44' '22.26'
Set-TextValue $ws 'E44' '  -0.28%  '
Set-TextValue $ws 'D45' '0.283'
Set-TextValue $ws 'E45' '  +19.98%  '
Set-TextValue $ws 'D46' '2.07'
Set-TextValue $ws 'E46' '  -0.61%  '
Set-TextValue $ws 'D47' '2.39'
Set-TextValue $ws 'E47' '  +2.74%  '
Set-TextValue $ws 'D48' '2.037.93'
Set-TextValue $ws 'E48' '  -1.05%  '
Set-TextValue $ws 'E49' '  +1.13%  '
Set-TextValue $ws 'D50' '0.0355'
Set-TextValue $ws 'E50' '  +11.91%  '
Set-TextValue $ws 'D51' '5.18'
Set-TextValue $ws 'E51' '  +2.68%  '

$excel.CutCopyMode = 0
